$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Add()
$ws.Name = "cat"
$ws.Range("A1").Value = "Outdoor Model"
$ws.Range("B1").Value = "Outdoor Quantity"
$ws.Range("C1").Value = "Outdoor Serial(s)"
$ws.Range("D1").Value = "Indoor Model"
$ws.Range("E1").Value = "Indoor Quantity"
$ws.Range("F1").Value = "Indoor Serial(s)"
